# Upload new version with timestamp
#
# A new shortage item ("سائل ريد") is inserted into the report table right
# before the existing "سرنجات 3 سم" row (its alphabetical slot in the
# Arabic item list). That pushes the two rows below it (items 44 and 45)
# down by one printed row, creates a brand-new row for what becomes item
# 46 (the data that used to belong to item 45's row), bumps the grand
# total by the new item's price, and refreshes the footer timestamp.
#
# Column A only ever holds the *printed row sequence number* (39, 40, 41,
# ...) - it is not tied to which item is shown - so a native Rows.Insert
# (which would shift those literal numbers down together with everything
# else) would not reproduce the target layout. Instead we rewrite each
# row's content explicitly, capturing old values before they are
# overwritten.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-TextValue($addr, $text) {
    # Forces the cell to be stored as a shared-string ("t=s") cell like the
    # target file expects, even though some of these columns carry a
    # numeric-looking NumberFormat. Restoring the original NumberFormat
    # afterwards keeps the cell's style bucket intact.
    $r = $ws.Range($addr)
    $fmt = $r.NumberFormat
    $r.NumberFormat = "@"
    $r.Value = $text
    $r.NumberFormat = $fmt
}

# ---------------------------------------------------------------------
# 1. Snapshot data that is about to move/disappear.
# ---------------------------------------------------------------------
$row49C = $ws.Range("C49").Text
$row49H = $ws.Range("H49").Text
$row49L = $ws.Range("L49").Text
$row49N = $ws.Range("N49").Text
$row49P = $ws.Range("P49").Text
$row49Q = $ws.Range("Q49").Text

$row50C = $ws.Range("C50").Text
$row50H = $ws.Range("H50").Text
$row50L = $ws.Range("L50").Text
$row50N = $ws.Range("N50").Text
$row50P = $ws.Range("P50").Text
$row50Q = $ws.Range("Q50").Text

$row51C = $ws.Range("C51").Text
$row51H = $ws.Range("H51").Text
$row51L = $ws.Range("L51").Text
$row51N = $ws.Range("N51").Text
$row51P = $ws.Range("P51").Text
$row51Q = $ws.Range("Q51").Text

$oldTotal = $ws.Range("P52").Value2
$oldFooterG = $ws.Range("G53").Text
$oldFooterK = $ws.Range("K53").Text

# ---------------------------------------------------------------------
# 2. Move the old total (row 52) and footer (row 53) rows down one slot
#    first, while row 52/53 still hold their original formatting.
# ---------------------------------------------------------------------
$ws.Range("P52:Q52").Copy()
$ws.Range("P53:Q53").PasteSpecial($xlPasteFormats)

$ws.Range("A53:K53").Copy()
$ws.Range("A54:K54").PasteSpecial($xlPasteFormats)

$excel.CutCopyMode = 0

# Drop the now-stale merges tied to the old row 52/53 layout before we
# repurpose those rows.
$ws.Range("P52:Q52").UnMerge()
$ws.Range("A53:F53").UnMerge()
$ws.Range("G53:I53").UnMerge()
$ws.Range("K53:Q53").UnMerge()

# ---------------------------------------------------------------------
# 3. Row 49 becomes the new item "سائل ريد".
# ---------------------------------------------------------------------
Set-TextValue "C49" "سائل ريد"
Set-TextValue "H49" "3:0"
Set-TextValue "L49" "0"
Set-TextValue "N49" "95.00"
Set-TextValue "P49" "95.0000"
Set-TextValue "Q49" "1:0"

# ---------------------------------------------------------------------
# 4. Rows 50 and 51 now carry what used to be in 49 and 50 (items 44 and
#    45's data); the item numbers printed in column A do not change.
# ---------------------------------------------------------------------
Set-TextValue "C50" $row49C
Set-TextValue "H50" $row49H
Set-TextValue "L50" $row49L
Set-TextValue "N50" $row49N
Set-TextValue "P50" $row49P
Set-TextValue "Q50" $row49Q

Set-TextValue "C51" $row50C
Set-TextValue "H51" $row50H
Set-TextValue "L51" $row50L
Set-TextValue "N51" $row50N
Set-TextValue "P51" $row50P
Set-TextValue "Q51" $row50Q

# ---------------------------------------------------------------------
# 5. Row 52 is a brand-new printed row (item 46), carrying what used to
#    be item 45's data (old row 51). Build it from row 51's current
#    per-column style pattern (identical styling to every item row).
# ---------------------------------------------------------------------
$ws.Range("A51:Q51").Copy()
$ws.Range("A52:Q52").PasteSpecial($xlPasteFormats)
$excel.CutCopyMode = 0
$ws.Rows.Item(52).RowHeight = $ws.Rows.Item(51).RowHeight

$ws.Range("A52").Value2 = 46
Set-TextValue "C52" $row51C
Set-TextValue "H52" $row51H
Set-TextValue "L52" $row51L
Set-TextValue "N52" $row51N
Set-TextValue "P52" $row51P
Set-TextValue "Q52" $row51Q

# ---------------------------------------------------------------------
# 6. Grand total row (now row 53): add the new item's price to the old
#    total; the row's own height changes too.
# ---------------------------------------------------------------------
$ws.Range("P53").Value2 = [double]$oldTotal + 95.0
$ws.Rows.Item(53).RowHeight = 24.75

# ---------------------------------------------------------------------
# 7. Footer row (now row 54): keep the page/author text, refresh the
#    generated timestamp.
# ---------------------------------------------------------------------
Set-TextValue "A54" "Sunday, 1 June, 2025 12:26 PM"
Set-TextValue "G54" $oldFooterG
Set-TextValue "K54" $oldFooterK
$ws.Rows.Item(54).RowHeight = 16.5

# ---------------------------------------------------------------------
# 8. Merge cells for the (re)built rows.
# ---------------------------------------------------------------------
$ws.Range("A52:B52").Merge()
$ws.Range("C52:G52").Merge()
$ws.Range("H52:K52").Merge()
$ws.Range("L52:M52").Merge()
$ws.Range("N52:O52").Merge()

$ws.Range("P53:Q53").Merge()

$ws.Range("A54:F54").Merge()
$ws.Range("G54:I54").Merge()
$ws.Range("K54:Q54").Merge()

$excel.CutCopyMode = 0
